$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "✉ john@company.com"
$ws.Range("C3").Value = "✉ jane@company.com"
$ws.Range("C4").Value = "✉ bob@company.com"
$ws.Range("C5").Value = "✉ alice@company.com"
